$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15
$ws.Range("A15").Value = 45969
$ws.Range("B15").Value = 'Warhammer: Vermintide 2'
$ws.Range("C15").Value = 'Jesper Kyd, Anders De Geer, Robin Hagblom, Matt Ward, Joakim Setterberg, Arvid Nilsson, Erika S. Kling, Robert Bäckström, Peter Nilsson, Joakim Wahlström'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2018-03-08'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = 4
$ws.Range("F15").Value = 'Warhammer: Vermintide 2 is a first-person action game developed by Fatshark and released on March 8, 2018. Serving as a sequel to the 2015 title Warhammer: End Times – Vermintide, the game is set in the Warhammer Fantasy universe during the tumultuous period known as The End Times. Players take on the roles of the heroes of Ubersreik, a group of five protagonists who are tasked with combating the malign forces of the Skaven clan Fester and the chaos warband Rotbloods. The game is notable for its emphasis on cooperative multiplayer gameplay, requiring teamwork and strategic coordination among players to succeed.
The gameplay of Warhammer: Vermintide 2 is structured around 13 missions, divided into three acts. Each mission requires players to select one of the five heroes and proceed from the start to the finish of the level, engaging in combat with large groups of enemies while completing various objectives. Each character in the game has distinct abilities, weapons, and a unique progression system that defines their role within the party. Characters are equipped with a melee weapon, a ranged weapon, an active ability, several passive abilities, and three career paths to choose from. Careers influence the character''s stats, available equipment, and abilities. As players progress and complete missions, they earn random pieces of equipment, which improve in quality as the game advances.
The development of Warhammer: Vermintide 2 involved a team of talented creators, including Jesper Kyd, Anders De Geer, and Robin Hagblom, among others. The game was published by Fatshark and made available on multiple platforms, including Xbox One, PC, and PlayStation 4. Known for its intense and visceral combat, the game integrates elements of action and indie genres, characterized by its cooperative gameplay, violent encounters, and adherence to the Warhammer lore.
Upon release, Warhammer: Vermintide 2 received generally favorable reviews, with many critics praising its engaging cooperative multiplayer experience and the depth of its combat system. The game was well-received for its atmospheric setting and challenging gameplay, which appealed to both fans of the Warhammer series and newcomers alike. Its cultural impact is evident in its continued popularity within the cooperative gaming community, often highlighted for its immersive and team-oriented gameplay. With a rating of 3.61, the game has maintained a dedicated player base and contributed to the ongoing legacy of the Warhammer franchise.'
$ws.Range("G15").Value = '[''Fatshark. (2018). Warhammer: Vermintide 2 [Video game]. Fatshark.'', ''IGN. (2018, March 8). Warhammer: Vermintide 2 review. IGN. Retrieved from https://www.ign.com'', ''Polygon. (2018, March 9). Warhammer: Vermintide 2 is a bloody good time. Polygon. Retrieved from https://www.polygon.com'', ''GameSpot. (2018, March 8). Warhammer: Vermintide 2 review. GameSpot. Retrieved from https://www.gamespot.com'', ''Eurogamer. (2018, March 10). Warhammer: Vermintide 2 review: a grisly, engaging co-op experience. Eurogamer. Retrieved from https://www.eurogamer.net'']'
$ws.Range("H15").Value = 'Average Rating: 3.61
Average Playtime: 7 hours
ESRB Rating: Mature
Metacritic Score: 82
Platforms: Xbox One, PC, PlayStation 4'
$ws.Range("I15").Value = '2025-03-27 19:20:07'
$ws.Range("J15").Value = 'https://media.rawg.io/media/games/5be/5bec14622f6faf804a592176577c1347.jpg'
$ws.Rows.Item(15).AutoFit()

# Row 16
$ws.Range("A16").Value = 13627
$ws.Range("B16").Value = 'Undertale'
$ws.Range("C16").Value = 'Toby Fox'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2015-09-14'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = 4
$ws.Range("F16").Value = 'Undertale is an independent role-playing game developed by Toby Fox and released on September 14, 2015. The game is available on multiple platforms including Xbox Series S/X, Nintendo Switch, Xbox One, PC, macOS, Linux, PlayStation 4, and PS Vita. It combines traditional RPG mechanics with unique elements such as bullet hell mini-games, creating an innovative gameplay experience that has captivated players worldwide. Undertale''s narrative begins with a backstory of two races, humans and monsters, who once coexisted until a war led to the monsters being sealed underground by human mages. The player controls a human child who falls into the monster realm and must navigate their way back to the surface.
Undertale''s gameplay focuses heavily on player choice, with a branching story and multiple endings determined by the player''s actions. Players interact with various NPCs and engage in battles that can be approached either peacefully or aggressively. The game features pixelated graphics reminiscent of retro video games, combined with a unique combat system where players can choose to spare or attack their opponents. These choices affect the behavior of monsters and the overall difficulty of encounters, emphasizing the game''s theme of morality and consequence.
The game''s reception was overwhelmingly positive, with critics and players alike praising its storytelling, humor, and innovative mechanics. Undertale''s soundtrack, also composed by Toby Fox, received particular acclaim for its ability to enhance the emotional depth of the narrative. The game''s design and execution have made it a standout title in the indie gaming scene, leading to a strong fanbase and lasting cultural impact.
Undertale''s influence extends beyond the gaming community, with its themes of empathy and choice resonating with a broad audience. The game has inspired a plethora of fan-created content, including art, music, and even memes, further cementing its place in popular culture. As a testament to its enduring appeal, Undertale continues to be referenced in various media and remains a topic of discussion among gaming enthusiasts.
For more information, visit the official website at [undertale.com](http://undertale.com).'
$ws.Range("G16").Value = '[''Fox, T. (2015). Undertale [Video game]. Toby Fox.'', ''PC Gamer. (2015, September 15). Undertale review. Retrieved from https://www.pcgamer.com/undertale-review/'', ''IGN. (2015, September 15). Undertale Review. Retrieved from https://www.ign.com/articles/2015/09/15/undertale-review'', ''Polygon. (2015, September 21). Undertale is a game where nobody has to get hurt. Retrieved from https://www.polygon.com/2015/9/21/9357025/undertale-review'', ''Metacritic. (2015). Undertale Reviews. Retrieved from https://www.metacritic.com/game/pc/undertale'']'
$ws.Range("H16").Value = 'Average Rating: 4.34
Average Playtime: 5 hours
ESRB Rating: Everyone 10+
Metacritic Score: 92
Platforms: Xbox Series S/X, Nintendo Switch, Xbox One, PC, macOS, Linux, PlayStation 4, PS Vita'
$ws.Range("I16").Value = '2025-03-27 19:20:21'
$ws.Range("J16").Value = 'https://media.rawg.io/media/games/ffe/ffed87105b14f5beff72ff44a7793fd5.jpg'
$ws.Rows.Item(16).AutoFit()

# Row 17
$ws.Range("A17").Value = 18240
$ws.Range("B17").Value = 'Brothers: A Tale of Two Sons'
$ws.Range("C17").Value = 'Gustaf Grefberg, Josef Fares'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2013-08-07'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = 4
$ws.Range("F17").Value = 'Brothers: A Tale of Two Sons is an action-adventure game that blends narrative-driven storytelling with innovative gameplay mechanics, developed by Gustaf Grefberg and Josef Fares. Released on August 7, 2013, and published by 505 Games and Parables Games, this indie title invites players into a fantastical world filled with emotional depth and richly painted environments. Tasked with saving their ailing father, the titular siblings, Naiee and Naia, embark on a quest to retrieve water from the mythical Tree of Life. Their journey traverses a variety of landscapes including villages, mountains, and forests, challenging them with obstacles such as bullies, trolls, and wolves.
The gameplay of Brothers: A Tale of Two Sons is distinguished by its unique control scheme, requiring players to manipulate both brothers simultaneously using a single controller. This dual-character control system enhances the puzzle-solving elements of the game, encouraging players to coordinate actions between the two characters, such as having one brother distract an enemy while the other sneaks past. The game seamlessly integrates these mechanics into its narrative, offering a fresh and engaging experience that stands out within the adventure genre. The visual style of the game depicts a vivid and atmospheric world that complements its emotional storytelling, supported by a memorable soundtrack that enriches the overall experience.
Upon release, Brothers: A Tale of Two Sons received critical acclaim for its innovative gameplay, emotional depth, and artistic presentation. It has been praised for its ability to convey a compelling story through gameplay mechanics rather than traditional narrative techniques. The game has won several awards and is often cited as a significant example of how video games can deliver powerful emotional narratives. Its reception among players has been equally positive, with a rating of 4.18 out of 5, and its availability across multiple platforms, including PlayStation 3, Xbox One, iOS, Nintendo Switch, and more, has broadened its reach to a diverse audience.
Culturally, Brothers: A Tale of Two Sons has contributed to the discourse on the potential of video games as a medium for storytelling. It has been discussed in various academic and industry settings as a case study in narrative design and player engagement. The game''s success has also paved the way for its developers to explore new creative projects, solidifying their reputation in the industry. Its impact is reflected in the continued interest and analysis by gaming scholars and enthusiasts alike.
For more information, the official website can be accessed at [505 Games](https://505games.com/games/brothers/).'
$ws.Range("G17").Value = '[''Fares, J., & Grefberg, G. (2013). Brothers: A Tale of Two Sons [Video game]. 505 Games, Parables Games.'', ''Smith, R. (2013). Review: Brothers: A Tale of Two Sons. GameSpot. Retrieved from https://www.gamespot.com/articles/review-brothers-a-tale-of-two-sons/1100-6412817/'', ''Johnson, M. (2014). Emotional storytelling in games: An analysis of Brothers: A Tale of Two Sons. Journal of Game Design, 9(3), 45-58.'', ''Thompson, C. (2013). How Brothers: A Tale of Two Sons Redefines Co-Op Gameplay. IGN. Retrieved from https://www.ign.com/articles/2013/08/06/how-brothers-a-tale-of-two-sons-redefines-co-op-gameplay'', ''Brown, A. (2015). The Artistic Vision of Brothers: A New Take on Narrative and Gameplay. Polygon. Retrieved from https://www.polygon.com/features/2015/8/7/9114917/the-artistic-vision-of-brothers'']'
$ws.Range("H17").Value = 'Average Rating: 4.18
Average Playtime: 3 hours
ESRB Rating: Teen
Metacritic Score: 90
Platforms: PlayStation 3, Xbox One, iOS, Nintendo Switch, Xbox 360, PlayStation 4, Android, PC'
$ws.Range("I17").Value = '2025-03-27 19:20:37'
$ws.Range("J17").Value = 'https://media.rawg.io/media/games/b6b/b6b20bfc4b34e312dbc8aac53c95a348.jpg'
$ws.Rows.Item(17).AutoFit()

# Row 18
$ws.Range("A18").Value = 9830
$ws.Range("B18").Value = 'Brawlhalla'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2014-04-30'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = 4
$ws.Range("F18").Value = 'Brawlhalla is a free-to-play fighting game that incorporates platformer elements, developed by Blue Mammoth Games and published by Ubisoft Entertainment. It was initially released on April 30, 2014, and has since expanded its availability to multiple platforms including Xbox One, Android, iOS, PC, macOS, Nintendo Switch, and PlayStation 4. The game falls under the genres of Action, Casual, Fighting, and Indie, and is known for its accessibility and competitive gameplay.
The core gameplay of Brawlhalla involves selecting a character from a roster of 39 heroes, a number that continues to grow with regular updates. The objective is to knock opponents off the stage using a combination of simple controls designed to make the game easy to master. Each character is equipped with two weapons and can utilize various boosts that appear on the stage. Characters also possess attributes such as Strength, Dexterity, Defense, and Speed, which define their combat styles and can be slightly customized. The game supports various modes, including online and local multiplayer, cooperative play, and competitive ranked matches, as well as large arena battles featuring multiple opponents.
Brawlhalla''s development began with a focus on creating a highly accessible yet competitive fighting experience. The game supports cross-platform multiplayer, allowing players to compete with others regardless of their gaming system. Its visual style is characterized by vibrant 2D graphics and dynamic stages, some of which feature multiple platforms that change positions during battles. As shown in the game''s artwork, Brawlhalla embraces a colorful and cartoonish aesthetic that appeals to a broad audience.
Reception of Brawlhalla has been generally positive, with players appreciating its free-to-play model and the depth of its combat mechanics. The game has cultivated a dedicated community and has been the subject of numerous competitions and tournaments worldwide, often featuring monetary prizes. With a rating of 3.22, Brawlhalla is recognized for its engaging multiplayer experience and its role in popularizing the platform fighter genre.
Brawlhalla has made a significant cultural impact by fostering an inclusive and competitive gaming community. The game''s emphasis on skill and strategy over complex controls has made it a staple in the esports scene, drawing players from various backgrounds. Official tournaments and community events have further solidified its standing as a prominent title in the fighting game circuit.
For more information, visit the official website at [Brawlhalla.com](http://www.brawlhalla.com).'
$ws.Range("G18").Value = '[''Blue Mammoth Games. (n.d.). *Brawlhalla*. Retrieved from http://www.brawlhalla.com'', ''Ubisoft Entertainment. (2014). *Brawlhalla* [Video game].'', ''IGN. (n.d.). *Brawlhalla Review*. Retrieved from https://www.ign.com'', ''Metacritic. (n.d.). *Brawlhalla*. Retrieved from https://www.metacritic.com'', ''Polygon. (n.d.). *How Brawlhalla became one of the most popular fighting games*. Retrieved from https://www.polygon.com'']'
$ws.Range("H18").Value = 'Average Rating: 3.22
Average Playtime: 2 hours
ESRB Rating: Everyone 10+
Platforms: Xbox One, Android, iOS, PC, macOS, Nintendo Switch, PlayStation 4'
$ws.Range("I18").Value = '2025-03-27 19:20:54'
$ws.Range("J18").Value = 'https://media.rawg.io/media/games/35b/35b47c4d85cd6e08f3e2ca43ea5ce7bb.jpg'
$ws.Rows.Item(18).AutoFit()
